$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "[공지] MBA AI/BigData 3rd term 문제 공개 – 비지니스와 AI 연계란?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ai-in-digital-marketing-example/#utm_source=rss&utm_medium=rss&utm_campaign=ai-in-digital-marketing-example"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D37").Value = "[Paper Review] data2vec: A General Framework for Self-supervised Learning in Speech, Vision and Language"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1935&mod=document&pageid=1"

$ws.Range("D51").Value = "[javascript] 오늘 날짜를 yyyy-MM-dd 형식으로 나타내기"
$ws.Range("E51").Value = "https://bskyvision.com/1243"

$ws.Range("D52").Value = "[R스터디:RQuestions] 첫 번째 시간/alluvial plot"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2546&utm_source=rss&utm_medium=rss&utm_campaign=r%25ec%258a%25a4%25ed%2584%25b0%25eb%2594%2594rquestions-%25ec%25b2%25ab-%25eb%25b2%2588%25ec%25a7%25b8-%25ec%258b%259c%25ea%25b0%2584-alluvial-plot"
